# SORTstatistics.xlsx - update processed data sheets
# - refresh raw measurement values on "ShinozakiSORT - Anlage" (cols B & D)
# - the cells that received new values also pick up the slightly larger
#   "Courier New" font (11pt instead of 10.5pt) used elsewhere in the sheet
# - dependent "Shinozaki - Anlage" sheet recalculates automatically (all of
#   its numbers are formulas off the first sheet)
# - move the active selection/tab to "Shinozaki - Anlage"

$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("ShinozakiSORT - Anlage")
$wsProc = $wb.Worksheets.Item("Shinozaki - Anlage")

# --- new raw values -------------------------------------------------------
$colB = @{
    4  = 14048
    5  = 20757
    6  = 24957
    7  = 28005
    8  = 30410
    9  = 32410
    10 = 34127
    11 = 35639
    12 = 36992
    13 = 38220
    14 = 39350
    15 = 40399
    16 = 41383
    17 = 42314
    18 = 43202
    19 = 44055
    20 = 44880
    21 = 45682
    22 = 46467
    23 = 47238
    24 = 48000
}

$colD = @{
    4  = 13667
    5  = 21152
    6  = 26409
    7  = 30416
    8  = 33616
    9  = 36256
    10 = 38485
    11 = 40398
    12 = 42064
    13 = 43530
    14 = 44833
    15 = 46000
}

foreach ($row in $colB.Keys) {
    $wsRaw.Range("B$row").Value = $colB[$row]
}

foreach ($row in $colD.Keys) {
    $wsRaw.Range("D$row").Value = $colD[$row]
}

# --- matching font bump on the touched columns -----------------------------
# Column C (untouched data) keeps its original 10.5pt Courier New style;
# columns B and D move to 11pt Courier New.
$wsRaw.Range("B4:B24").Font.Size = 11
$wsRaw.Range("D4:D15").Font.Size = 11

# --- selection / active sheet ---------------------------------------------
$wsRaw.Range("D4").Select()

$wsProc.Activate()
$wsProc.Range("F13").Select()

$wb.RefreshAll()
